$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.504.10'
$ws.Range('E2').Value = '  +2.32%  '

$ws.Range('D3').Value = '2.697.63'
$ws.Range('E3').Value = '  +2.52%  '

$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '525.41'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.55%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '145.07'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -1.00%  '

$ws.Range('E7').Value = '  +0.09%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.577'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +1.37%  '

$ws.Range('D9').Value = '2.722.35'
$ws.Range('E9').Value = '  +2.49%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.69'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +5.90%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.105'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.01%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.339'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.15%  '

$ws.Range('E13').Value = '  +3.06%  '

$ws.Range('D14').Value = '3.176.01'
$ws.Range('E14').Value = '  +2.46%  '

$ws.Range('D15').Value = '60.535.75'
$ws.Range('E15').Value = '  +2.44%  '

$ws.Range('B16').Value = 'Avalanche'
$ws.Range('C16').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '21.26'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +1.12%  '

$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '2.796.62'
$ws.Range('E17').Value = '  +5.45%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.0000137'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.32%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '345.51'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -1.10%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.51'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.29%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '10.63'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +2.77%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.45'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +4.29%  '

$ws.Range('E23').Value = '  +0.03%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '63.55'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +2.93%  '

$ws.Range('E25').Value = '  +0.43%  '

$ws.Range('E26').Value = '  +4.00%  '

$ws.Range('E27').Value = '  -0.15%  '

$ws.Range('D28').Value = '0.0₃0818'
$ws.Range('E28').Value = '  +1.01%  '

$ws.Range('E30').Value = '  +8.35%  '

$ws.Range('E32').Value = '  +0.57%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '19.03'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.06%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '149.96'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.07%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.24'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +5.74%  '

$ws.Range('E36').Value = '  +8.34%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.945'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -3.75%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.52'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +7.50%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.870'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +2.69%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '37.11'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.80%  '

$ws.Range('E41').Value = '  -0.72%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '282.56'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +1.18%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '20.11'
$ws.Range('D43').Style = "Normal"

$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.996'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.19%  '

$ws.Range('B45').Value = 'Stellar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0987'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.02%  '

$ws.Range('D46').Value = '2.139.51'
$ws.Range('E46').Value = '  +7.47%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.609'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.37%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0539'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +2.53%  '

$ws.Range('B49').Value = 'WhiteBITCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '10.47'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +1.88%  '

$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '4.79'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.62%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0232'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.86%  '
